$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.185.18"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.943.73"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'376.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'102.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").Value = "'36.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "'0.0838"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "3.404.18"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'7.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "2.928.80"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "'0.972"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "51.105.01"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("E19").Value = "  -6.62%  "
$ws.Range("D20").Value = "'7.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("D21").Value = "'12.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'263.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'68.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").Value = "'2.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "'7.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.78%  "
$ws.Range("D27").Value = "'8.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.22%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.168"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.113"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.16%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'25.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'9.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'34.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").Value = "'0.0458"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "'50.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'16.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "'121.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'21.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "'0.274"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "'3.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "1.999.62"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").Value = "'0.0349"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "'5.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.24%  "
